$paraFragments = @(
    '',
    '<w:r><w:t>$(document).ready(function () {</w:t></w:r>',
    '<w:r><w:t xml:space="preserve">    // Attach click event to all elements</w:t></w:r>',
    '<w:r><w:t xml:space="preserve">    $(''body'').on(''click'', function (e) {</w:t></w:r>',
    '<w:r><w:t xml:space="preserve">        // Check if __</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>doPostBack</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> function exists and the click event target is not null</w:t></w:r>',
    '<w:r><w:t xml:space="preserve">        if (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>typeof</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> __</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>doPostBack</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> === ''function'' &amp;&amp; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>e.target</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> !== null) {</w:t></w:r>',
    '<w:r><w:t xml:space="preserve">            // Get the ID of the clicked element</w:t></w:r>',
    '<w:r><w:t xml:space="preserve">            </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>var</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>elementId</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>e.target.id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>;</w:t></w:r>',
    '',
    '<w:r><w:t xml:space="preserve">            // Check if the clicked element caused a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>postback</w:t></w:r><w:proofErr w:type="spellEnd"/>',
    '<w:r><w:t xml:space="preserve">            if (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>elementId</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> !== '''') {</w:t></w:r>',
    '<w:r><w:t xml:space="preserve">                // Log the ID of the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>postback</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> element</w:t></w:r>',
    '<w:r><w:t xml:space="preserve">                </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>console.log</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(''</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Postback</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> element ID: '' + </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>elementId</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>);</w:t></w:r>',
    '<w:r><w:t xml:space="preserve">            }</w:t></w:r>',
    '<w:r><w:t xml:space="preserve">        }</w:t></w:r>',
    '<w:r><w:t xml:space="preserve">    });</w:t></w:r>',
    '<w:r><w:t>});</w:t></w:r>'
)

$d = $word.ActiveDocument

function Insert-ParagraphXml($range, $innerXml) {
    $xmlFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $innerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($xmlFrag)
}

# Start right after the last paragraph ("</html>") and create one working
# empty paragraph to anchor each subsequent insertion.
$tail = $d.Paragraphs.Last.Range
$tail.Collapse(0)
$tail.InsertParagraphAfter()

$cur = $d.Paragraphs.Last.Range
$cur.Collapse(0)

foreach ($frag in $paraFragments) {
    Insert-ParagraphXml $cur $frag
    $cur = $d.Paragraphs.Last.Range
    $cur.Collapse(0)
}

# Remove the now-superfluous trailing empty paragraph created as the
# anchor for the final insertion.
$trailing = $d.Paragraphs.Last.Range
$delRange = $d.Range($trailing.Start - 1, $trailing.End)
$delRange.Delete()
